$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" '27.097.45'
Set-TextCell "E2" '  +3.04%  '
Set-TextCell "D3" '1.657.80'
Set-TextCell "E3" '  +3.53%  '
Set-TextCell "D4" '0.998'
Set-TextCell "E4" '  -0.25%  '
Set-TextCell "D5" '215.63'
Set-TextCell "E5" '  +1.41%  '
Set-TextCell "D7" '0.998'
Set-TextCell "E7" '  -0.28%  '
Set-TextCell "E8" '  +2.12%  '
Set-TextCell "E9" '  +1.19%  '
Set-TextCell "D10" '19.52'
Set-TextCell "E10" '  +2.69%  '
Set-TextCell "D11" '0.0860'
Set-TextCell "D12" '1.888.94'
Set-TextCell "E12" '  +3.41%  '
Set-TextCell "D13" '1.654.80'
Set-TextCell "E13" '  +4.98%  '
Set-TextCell "D14" '4.08'
Set-TextCell "E15" '  +2.55%  '
Set-TextCell "D16" '64.99'
Set-TextCell "E16" '  +1.99%  '
Set-TextCell "D17" '241.81'
Set-TextCell "E17" '  +4.97%  '
Set-TextCell "D18" '27.091.37'
Set-TextCell "E18" '  +3.02%  '
Set-TextCell "D19" '7.87'
Set-TextCell "E19" '  +3.29%  '
Set-TextCell "D20" '0.0₃0729'
Set-TextCell "E20" '  +1.08%  '
Set-TextCell "D21" '0.997'
Set-TextCell "E21" '  -0.36%  '
Set-TextCell "D22" '4.45'
Set-TextCell "E22" '  +3.88%  '
Set-TextCell "E23" '  +5.32%  '
Set-TextCell "E24" '  +3.40%  '
Set-TextCell "D25" '146.29'
Set-TextCell "E25" '  +0.40%  '
Set-TextCell "D26" '0.999'
Set-TextCell "E26" '  -0.25%  '
Set-TextCell "E27" '  +2.37%  '
Set-TextCell "E28" '  +1.19%  '
Set-TextCell "D29" '15.86'
Set-TextCell "E29" '  +2.75%  '
Set-TextCell "E30" '  +0.60%  '
Set-TextCell "E31" '  +0.76%  '
Set-TextCell "D32" '1.522.40'
Set-TextCell "E32" '  +5.15%  '
Set-TextCell "E33" '  +2.55%  '
Set-TextCell "E34" '  +2.81%  '
Set-TextCell "D35" '1.56'
Set-TextCell "E35" '  +6.20%  '
Set-TextCell "D36" '2.42'
Set-TextCell "E36" '  -0.18%  '
Set-TextCell "D37" '0.579'
Set-TextCell "E37" '  +1.48%  '
Set-TextCell "E38" '  +7.88%  '
Set-TextCell "E39" '  +2.24%  '
Set-TextCell "D40" '5.97'
Set-TextCell "E40" '  +3.16%  '
Set-TextCell "E41" '  -0.37%  '
Set-TextCell "D42" '2.33'
Set-TextCell "E42" '  +6.91%  '
Set-TextCell "D43" '64.94'
Set-TextCell "E43" '  +6.92%  '
Set-TextCell "D44" '1.797.16'
Set-TextCell "E44" '  +3.31%  '
Set-TextCell "E45" '  +1.98%  '
Set-TextCell "D46" '0.914'
Set-TextCell "E46" '  -0.82%  '
Set-TextCell "D47" '90.52'
Set-TextCell "E47" '  +3.60%  '
Set-TextCell "E48" '  +3.47%  '
Set-TextCell "D49" '0.0981'
Set-TextCell "E49" '  +3.21%  '
Set-TextCell "D51" '7.54'
Set-TextCell "E51" '  +1.68%  '
